$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet (tab) to the English name
$ws.Name = "DEFAULT LIST"

# 2. Translate the table header row to English.
#    (Updating the cell values also keeps the Tabela1 table column names in sync.)
$ws.Range("B1").Value = "Name (Optional)"
$ws.Range("C1").Value = "Telephone"
$ws.Range("D1").Value = "Sent"

# 3. Add the four "reminder" notes to the right of the table (columns F:I),
#    each merged across F:I, bold and horizontally centered.
$ws.Range("F2").Value = "Sempre salve e feche antes de prosseguir."
$ws.Range("F3").Value = "Para usar mais abas, copie esta."
$ws.Range("F5").Value = "Always save and close before proceeding."
$ws.Range("F6").Value = "To use more tabs, copy this one."

$notesTop = $excel.Union($ws.Range("F2:I2"), $ws.Range("F3:I3"))
$notesBottom = $excel.Union($ws.Range("F5:I5"), $ws.Range("F6:I6"))

$notesTop.HorizontalAlignment = -4108
$notesTop.Font.Bold = $true
$notesBottom.HorizontalAlignment = -4108
$notesBottom.Font.Bold = $true

$ws.Range("F2:I2").Merge() | Out-Null
$ws.Range("F3:I3").Merge() | Out-Null
$ws.Range("F5:I5").Merge() | Out-Null
$ws.Range("F6:I6").Merge() | Out-Null

# 4. Update the active selection to match the edited workbook.
$ws.Range("D2:D6").Select() | Out-Null
